$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.955.73"
$ws.Range("E2").Value = "  +0.30%  "

$ws.Range("D3").Value = "1.884.93"
$ws.Range("E3").Value = "  +1.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5144"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3739"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07191"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8982"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07653"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.35%  "

$ws.Range("D13").Value = "1.868.97"
$ws.Range("E13").Value = "  +0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("E15").Value = "  +0.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008479"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").Value = "26.996.61"
$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.035"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").Value = "2.127.01"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.381"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.289"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.726"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("E30").Value = "  +4.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.771"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.88%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09172"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.228"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7645"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.34%  "

$ws.Range("E36").Value = "  +0.99%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.580"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5570"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01982"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.99%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.012"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.602"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1497"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4805"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.0000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.598"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.53%  "

$ws.Range("E51").Value = "  +1.33%  "

$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"

